# Tipuri de date pentru fiecare coloana
# Update the data type cells for "day/month/year"-style columns from
# INT(2) UNSIGNED / INT(4) UNSIGNED to CHAR(2) / CHAR(4) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Medici 2NF table (row 11): zin_medic / lunan_medic / ann_medic types
$ws.Range("D11").Value = "CHAR(2)"
$ws.Range("E11").Value = "CHAR(2)"
$ws.Range("F11").Value = "CHAR(4)"

# ContracteMedici table (row 15): ziua_angajarii / luna_angajarii / anul_angajarii types
$ws.Range("D15").Value = "CHAR(2)"
$ws.Range("E15").Value = "CHAR(2)"
$ws.Range("F15").Value = "CHAR(4)"

# Programari_Consultatii 2NF table (row 27): zi_programare / luna_programare / an_programare /
# ora_programare / minut_programare types
$ws.Range("B27").Value = "CHAR(2)"
$ws.Range("C27").Value = "CHAR(2)"
$ws.Range("D27").Value = "CHAR(4)"
$ws.Range("E27").Value = "CHAR(2)"
$ws.Range("F27").Value = "CHAR(2)"

# Fise_Pacienti 2NF table (row 39): ziua_nasterii / luna_nasterii / anul_nasterii types
$ws.Range("D39").Value = "CHAR(2)"
$ws.Range("E39").Value = "CHAR(2)"
$ws.Range("F39").Value = "CHAR(4)"

# Update the stored selection to match the authored workbook state
$ws.Range("D45").Select()
